# Commit: "Fixed Tests for SamplePatholoy, SampleType, Se, StageOfDisease, and Study"
#
# The "Cases" Cypher query stored in cell B2 of the "startup" sheet is
# trimmed: the last RETURN line that projected the cohort description
# (`coalesce(co.cohort_description, '') AS `Cohort``) is removed, along
# with the now-dangling trailing comma on the preceding line.
#
# (Everything else in the underlying OOXML diff -- shared-string reordering,
# row-height/column-width recalculation, Excel build/version stamps, GUIDs,
# and window geometry -- are save-artifacts of the authoring Excel build and
# fall out automatically once the cell content below is written.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newCasesQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n" +
                 "MATCH (c)<--(diag:diagnosis)`n" +
                 "MATCH (samp:sample)-->(c) `n" +
                 "WHERE samp.summarized_sample_type IN [`"Whole Blood`"] `n" +
                 "OPTIONAL MATCH (co:cohort)<-[*]-(c)`n" +
                 "  WITH DISTINCT c, s, demo, diag, co`n" +
                 "RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n" +
                 "        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n" +
                 "        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n" +
                 "        coalesce(demo.breed, '') AS Breed ,`n" +
                 "        coalesce(diag.disease_term, '') AS Diagnosis ,`n" +
                 "        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n" +
                 "        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n" +
                 "        coalesce(demo.sex, '') AS Sex ,`n" +
                 "        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n" +
                 "        coalesce(demo.weight, '') AS ``Weight (kg)``,`n" +
                 "        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $newCasesQuery

# Matches the author's final cursor position recorded in the worksheet's
# <selection> after editing the formula in B2.
$ws.Range("B2").Select()
